# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Fri Sep  1 10:35:11 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.139.79"
$ws.Range("E2").Value = "  -4.37%  "
$ws.Range("D3").Value = "1.650.70"
$ws.Range("E3").Value = "  -3.54%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5108"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.40%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06433"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07784"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").Value = "1.656.54"
$ws.Range("E12").Value = "  -3.19%  "
$ws.Range("E13").Value = "  -4.93%  "
$ws.Range("D14").Value = "1.878.40"
$ws.Range("E14").Value = "  -3.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5515"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.78%  "
$ws.Range("D16").Value = "0.0₅8004"
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.77%  "
$ws.Range("D18").Value = "26.151.12"
$ws.Range("E18").Value = "  -4.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.387"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.38%  "
$ws.Range("E22").Value = "  -3.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.036"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.751"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1175"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.58%  "
$ws.Range("E28").Value = "  -3.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05133"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.83%  "
$ws.Range("E31").Value = "  -3.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.346"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.212"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.558"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.739"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9232"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.351"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Value = "1.169.00"
$ws.Range("E38").Value = "  +1.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5687"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01584"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.552"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.655"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8227"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("D46").Value = "1.788.63"
$ws.Range("E46").Value = "  -3.50%  "
$ws.Range("D47").Value = "0.0₈116"
$ws.Range("E47").Value = "  +1.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4552"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.44"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.78%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.005"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  -3.02%  "
